# Add columns I (I0) and J (IF) following the existing header row's formatting
# (bold, centered, top-aligned, thin-bordered) by copying H1's format, then
# fill in the header labels and the two data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header labels
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the header style (font/alignment/border) from the existing H1 header
# cell onto the two new header cells so they reuse the same cell style.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data rows
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 5

$ws.Range("I3").Value = 7
$ws.Range("J3").Value = 8
